# Apply updated cryptocurrency price/volume data per the Jun 29 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.753.72"
$ws.Range("E2").Value = "  -1.39%  "

# Row 3
$ws.Range("D3").Value = "3.383.82"
$ws.Range("E3").Value = "  -1.96%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").Value = "'569.17"
$ws.Range("E5").Value = "  -1.80%  "

# Row 6
$ws.Range("D6").Value = "'140.74"
$ws.Range("E6").Value = "  -2.68%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "3.383.83"
$ws.Range("E8").Value = "  -1.90%  "

# Row 9
$ws.Range("E9").Value = "  -0.57%  "

# Row 10
$ws.Range("D10").Value = "'7.49"
$ws.Range("E10").Value = "  -1.39%  "

# Row 11
$ws.Range("E11").Value = "  -1.59%  "

# Row 12
$ws.Range("D12").Value = "'0.395"
$ws.Range("E12").Value = "  +1.87%  "

# Row 13
$ws.Range("D13").Value = "3.961.56"
$ws.Range("E13").Value = "  -2.01%  "

# Row 14
$ws.Range("D14").Value = "'28.53"
$ws.Range("E14").Value = "  +1.57%  "

# Row 15
$ws.Range("E15").Value = "  +2.29%  "

# Row 16
$ws.Range("D16").Value = "'0.0000170"
$ws.Range("E16").Value = "  -1.55%  "

# Row 17
$ws.Range("D17").Value = "3.383.72"
$ws.Range("E17").Value = "  -2.00%  "

# Row 18
$ws.Range("D18").Value = "60.831.58"
$ws.Range("E18").Value = "  -1.46%  "

# Row 19
$ws.Range("E19").Value = "  -0.85%  "

# Row 20
$ws.Range("D20").Value = "'13.96"
$ws.Range("E20").Value = "  -1.98%  "

# Row 21
$ws.Range("D21").Value = "'9.00"
$ws.Range("E21").Value = "  -5.46%  "

# Row 22
$ws.Range("D22").Value = "'383.44"
$ws.Range("E22").Value = "  -1.62%  "

# Row 23
$ws.Range("D23").Value = "'0.559"
$ws.Range("E23").Value = "  -0.85%  "

# Row 24
$ws.Range("D24").Value = "'73.62"
$ws.Range("E24").Value = "  +0.28%  "

# Row 25
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.38%  "

# Row 26
$ws.Range("D26").Value = "'0.0000117"
$ws.Range("E26").Value = "  -5.40%  "

# Row 27
$ws.Range("D27").Value = "3.518.55"
$ws.Range("E27").Value = "  -1.93%  "

# Row 28
$ws.Range("E28").Value = "  -0.45%  "

# Row 29
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.09%  "

# Row 30
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'7.39"
$ws.Range("E30").Value = "  -3.23%  "

# Row 31
$ws.Range("D31").Value = "'8.00"
$ws.Range("E31").Value = "  -1.82%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.14"
$ws.Range("E32").Value = "  -1.86%  "

# Row 33
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.43"
$ws.Range("E33").Value = "  -2.96%  "

# Row 34
$ws.Range("E34").Value = "  -0.01%  "

# Row 35
$ws.Range("D35").Value = "'23.58"
$ws.Range("E35").Value = "  -1.93%  "

# Row 36
$ws.Range("D36").Value = "'6.96"
$ws.Range("E36").Value = "  -0.43%  "

# Row 37
$ws.Range("D37").Value = "'166.16"
$ws.Range("E37").Value = "  -0.46%  "

# Row 38
$ws.Range("D38").Value = "3.413.86"
$ws.Range("E38").Value = "  -1.92%  "

# Row 39
$ws.Range("D39").Value = "'4.98"
$ws.Range("E39").Value = "  -2.97%  "

# Row 40
$ws.Range("E40").Value = "  -4.42%  "

# Row 41
$ws.Range("D41").Value = "'27.86"
$ws.Range("E41").Value = "  -0.93%  "

# Row 42
$ws.Range("D42").Value = "'0.0773"
$ws.Range("E42").Value = "  -0.99%  "

# Row 43
$ws.Range("E43").Value = "  +0.05%  "

# Row 44
$ws.Range("E44").Value = "  -3.05%  "

# Row 45
$ws.Range("D45").Value = "'41.88"
$ws.Range("E45").Value = "  -1.01%  "

# Row 46
$ws.Range("E46").Value = "  -1.49%  "

# Row 47
$ws.Range("E47").Value = "  -3.59%  "

# Row 48
$ws.Range("E48").Value = "  -2.67%  "

# Row 49
$ws.Range("D49").Value = "2.499.68"
$ws.Range("E49").Value = "  -3.45%  "

# Row 50
$ws.Range("D50").Value = "'23.51"
$ws.Range("E50").Value = "  +2.27%  "

# Row 51
$ws.Range("E51").Value = "  -1.29%  "
